# Update the "user.followers_count" description cell (D12) in the data
# dictionary sheet: reword the text and make "esta siguiendo" bold, to
# match the corrected description of the field.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$cell = $ws.Range("D12")

$part1 = "El número de usuarios que esta cuenta "
$part2 = "esta siguiendo"
$part3 = " (también conocido como sus ""seguidores""). Bajo ciertas condiciones de coacción, este campo indicará temporalmente ""0""."

# Set the full plain text first.
$cell.Value = "$part1$part2$part3"

$len1 = $part1.Length
$len2 = $part2.Length
$len3 = $part3.Length

# Make the middle run ("esta siguiendo") bold, keeping the same font
# used by the rest of the cell.
$boldRun = $cell.Characters($len1 + 1, $len2)
$boldRun.Font.Name = "Calibri"
$boldRun.Font.Size = 12
$boldRun.Font.Bold = $true
$boldRun.Font.Color = 2367776   # RGB(0x20,0x21,0x24) -> FF202124

# Re-assert the regular formatting for the trailing run so every run in
# the cell carries explicit, matching font information.
$tailRun = $cell.Characters($len1 + $len2 + 1, $len3)
$tailRun.Font.Name = "Calibri"
$tailRun.Font.Size = 12
$tailRun.Font.Bold = $false
$tailRun.Font.Color = 2367776

# Restore the view state: scrolled down so row 7 is at the top, with
# F11 as the active/selected cell.
$ws.Activate() | Out-Null
$ws.Range("F11").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
